# Applies the "CU Consultar proximos pagos" update to the
# "Plantilla de Casos de Uso" workbook.
#
# Content changes (sheet "Casos de Uso"):
#   - Row 8 ("Consultar proximos pagos de alumnos"):
#       Estado (E8):   "planificado" -> "Hecho"
#       Esfuerzo (F8): 1 -> 5
#       Esfuerzo% (G8): 0 -> 100
#   - Row 7 ("Generar recibo de pago"):
#       Esfuerzo% (G7): 0 -> 90
#   - Selected/active cell on the sheet becomes H12 (was F14)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")
$ws.Activate()

# Row 8: CU "Consultar proximos pagos de alumnos"
$ws.Range("E8").Value = "Hecho"
$ws.Range("F8").Value = 5
$ws.Range("G8").Value = 100

# Row 7: CU "Generar recibo de pago"
$ws.Range("G7").Value = 90

# Update the selected cell to match the saved view state
$ws.Range("H12").Select()

$wb.Save()
